# Refresh the crypto price/volume listing (rows 2-51) to the latest
# scraped snapshot. A couple of rows also swapped rank order
# (WrappedEther/Polkadot and Stellar/ARBITRUM), so those rows' Coin,
# Link, Price and Volume columns are fully rewritten rather than just
# the numbers.
#
# NumberFormat is forced to text ("@") before writing any Price value
# that would otherwise be auto-parsed by Excel as a number (e.g.
# "1.009", "6.527"), so it is stored as the exact original string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.525.35'
$ws.Range('E2').Value = '  -3.50%  '
$ws.Range('D3').Value = '1.997.94'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.009'
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '329.85'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5007'
$ws.Range('E7').Value = '  -4.61%  '
$ws.Range('E8').Value = '  -5.30%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '53.49'
$ws.Range('E9').Value = '  -2.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08990'
$ws.Range('E10').Value = '  -4.68%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.123'
$ws.Range('E11').Value = '  -5.17%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '23.30'
$ws.Range('E12').Value = '  -8.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.087'
$ws.Range('E13').Value = '  -7.58%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.527'
$ws.Range('E14').Value = '  -6.63%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '1.946.95'
$ws.Range('E15').Value = '  -8.18%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '95.79'
$ws.Range('E16').Value = '  -6.36%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.008'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001118'
$ws.Range('E18').Value = '  -4.58%  '
$ws.Range('E19').Value = '  -1.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.74'
$ws.Range('E20').Value = '  -8.57%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.008'
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.967'
$ws.Range('D23').Value = '29.527.32'
$ws.Range('E23').Value = '  -3.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.94'
$ws.Range('E24').Value = '  -6.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.270'
$ws.Range('E25').Value = '  -2.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.78'
$ws.Range('E26').Value = '  -2.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.68'
$ws.Range('E27').Value = '  -7.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.566'
$ws.Range('E28').Value = '  -4.76%  '
$ws.Range('E29').Value = '  -8.94%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '127.81'
$ws.Range('E30').Value = '  -5.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.051'
$ws.Range('E31').Value = '  -9.84%  '
$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.578'
$ws.Range('E32').Value = '  -11.95%  '
$ws.Range('B33').Value = 'Stellar'
$ws.Range('C33').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09945'
$ws.Range('E33').Value = '  -6.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.851'
$ws.Range('E34').Value = '  -7.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.783'
$ws.Range('E35').Value = '  -4.56%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.609'
$ws.Range('E36').Value = '  -9.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02463'
$ws.Range('E37').Value = '  -7.50%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06355'
$ws.Range('E38').Value = '  -7.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.293'
$ws.Range('E39').Value = '  -3.46%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6537'
$ws.Range('E40').Value = '  -8.80%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.67'
$ws.Range('E41').Value = '  -8.26%  '
$ws.Range('E42').Value = '  -8.43%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6333'
$ws.Range('E44').Value = '  -8.89%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.210'
$ws.Range('E45').Value = '  -7.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.271'
$ws.Range('E47').Value = '  -6.28%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.529'
$ws.Range('E48').Value = '  -3.63%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07000'
$ws.Range('E49').Value = '  -3.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00000000322'
$ws.Range('E50').Value = '  -7.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.136'
$ws.Range('E51').Value = '  -6.15%  '
